# Add team record (Wins/Losses/Ties) columns to the data sheet.
# Mirrors the commit: "Added team record to data" - W/L/T appended as
# new columns (AD, AE, AF) on the same sheet, repeated on every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the existing header style (same as A1:AC1)
# onto the three new header cells, then set their labels.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-52): every row gets the same team record.
$wins = 95
$losses = 68
$ties = 0

for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
